# "fix the multigrid styling issues"
#
# Adds a fixed "employee id" style date column (col C) alongside the
# existing Start/End date column (col B), fills in the new progress rows
# (16-20) with their Start/End dates, widens column B slightly to fit the
# new dates, and leaves the cursor on the newly-entered C19 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = "d-mmm-yy"

# --- Row 16: "commit code" gains Start/End dates ------------------------
$ws.Range("B16").Value = 45214
$ws.Range("B16").NumberFormat = $dateFmt
$ws.Range("C16").Value = 45214
$ws.Range("C16").NumberFormat = $dateFmt

# --- Row 17: "integrate with List Employee screen" gains Start/End dates
$ws.Range("B17").Value = 45214
$ws.Range("B17").NumberFormat = $dateFmt
$ws.Range("C17").Value = 45214
$ws.Range("C17").NumberFormat = $dateFmt

# --- Row 18 (new): multigrid work item -----------------------------------
$ws.Range("A18").Value = "use multi grid with fixed employee id column to display the data"
$ws.Range("B18").Value = 45214
$ws.Range("B18").NumberFormat = $dateFmt
$ws.Range("C18").Value = 45224
$ws.Range("C18").NumberFormat = $dateFmt

# --- Row 19 (new): styling fix work item ---------------------------------
$ws.Range("A19").Value = "fix the styling issues "
$ws.Range("B19").Value = 45240
$ws.Range("B19").NumberFormat = $dateFmt
$ws.Range("C19").Value = 45240
$ws.Range("C19").NumberFormat = $dateFmt

# --- Row 20 (new): testing work item (no dates yet) ----------------------
$ws.Range("A20").Value = "testing the screen "

# Column B needs to widen slightly (10 chars) to fit the new dates; column C
# keeps the existing width used by the original Start/End date column.
$ws.Columns.Item(2).ColumnWidth = 9.14

# Final cursor position/selection left on the newly entered C19 cell.
$ws.Range("C19").Select()
